$wb = $excel.ActiveWorkbook

# --- Sheet "Canada": add row 13 (2020-12-01 / 44166) ---
$ws1 = $wb.Worksheets.Item("Canada")
$ws1.Range("A13").NumberFormat = "d-mmm-yy"
$ws1.Range("A13").Value = 44166
$ws1.Range("B13").NumberFormat = "d-mmm-yy"
$ws1.Range("B13").Value = "Canada"
$ws1.Range("C13").Value = 53.6
$ws1.Range("D13").Value = 1755.8
$ws1.Range("C14").Select()

# --- Sheet "Province": add rows 112-121 (2020-12-01 / 44166) ---
$ws2 = $wb.Worksheets.Item("Province")

$provinceData = @(
    @("Newfoundland & Labrador", 6.1, 31.4),
    @("Prince Edward Island", 24.6, 8.6),
    @("Nova Scotia", 9.3000000000000007, 43.4),
    @("New Brunswick", 26.8, 36.4),
    @("Quebec", 26.9, 305.8),
    @("Ontario", 80, 762.5),
    @("Manitoba", 61.9, 55.7),
    @("Saskatchewan", 33.5, 47),
    @("Alberta", 53.9, 271.39999999999998),
    @("British Columbia", 51.3, 193.7)
)

$row = 112
$first = $true
foreach ($entry in $provinceData) {
    $ws2.Range("A$row").NumberFormat = "d-mmm-yy"
    $ws2.Range("A$row").Value = 44166
    if ($first) {
        $ws2.Range("B$row").NumberFormat = "d-mmm-yy"
    }
    $ws2.Range("B$row").Value = $entry[0]
    $ws2.Range("C$row").Value = $entry[1]
    $ws2.Range("D$row").Value = $entry[2]
    $first = $false
    $row = $row + 1
}

$ws2.Application.ActiveWindow.ScrollRow = 98
$ws2.Range("C122").Select()
